# Apply targeted numeric cell updates across multiple worksheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 28468.805
$ws.Range("J17").Value = 28468.805
$ws.Range("L17").Value = 85406.41500000001
$ws.Range("N17").Value = -85742.41500000001
$ws.Range("H19").Value = 10989706
$ws.Range("I19").Value = 17857744
$ws.Range("J19").Value = 846
$ws.Range("K19").Value = 17857744
$ws.Range("L19").Value = 846
$ws.Range("M19").Value = -17857569
$ws.Range("N19").Value = -1196
$ws.Range("H31").Value = 1280
$ws.Range("I31").Value = 800
$ws.Range("J31").Value = 1600
$ws.Range("K31").Value = 2400
$ws.Range("L31").Value = 4800
$ws.Range("M31").Value = -2170
$ws.Range("N31").Value = -5260
$ws.Range("H38").Value = 1509.6
$ws.Range("I38").Value = 220.5
$ws.Range("J38").Value = 2982.8572
$ws.Range("K38").Value = 661.5
$ws.Range("L38").Value = 8948.571599999999
$ws.Range("M38").Value = -289.5
$ws.Range("N38").Value = -9692.571599999999
$ws.Range("H39").Value = 1381.875
$ws.Range("I39").Value = 74.333336
$ws.Range("J39").Value = 2166.4
$ws.Range("K39").Value = 223.000008
$ws.Range("L39").Value = 6499.200000000001
$ws.Range("M39").Value = 72.99999199999999
$ws.Range("N39").Value = -7091.200000000001
$ws.Range("H70").Value = 1559.8206
$ws.Range("I70").Value = 950.2174
$ws.Range("J70").Value = 2436.125
$ws.Range("K70").Value = 2850.6522
$ws.Range("L70").Value = 7308.375
$ws.Range("M70").Value = -2580.6522
$ws.Range("N70").Value = -7848.375
$ws.Range("H73").Value = 1559.8206
$ws.Range("I73").Value = 950.2174
$ws.Range("J73").Value = 2436.125
$ws.Range("K73").Value = 2850.6522
$ws.Range("L73").Value = 7308.375
$ws.Range("M73").Value = -1914.6522
$ws.Range("N73").Value = -9180.375
$ws.Range("H80").Value = 3684.1924
$ws.Range("I80").Value = 381.18182
$ws.Range("J80").Value = 6106.4
$ws.Range("K80").Value = 1143.54546
$ws.Range("L80").Value = 18319.2
$ws.Range("M80").Value = -145.54546
$ws.Range("N80").Value = -20315.2
$ws.Range("H83").Value = 3684.1924
$ws.Range("I83").Value = 381.18182
$ws.Range("J83").Value = 6106.4
$ws.Range("K83").Value = 3430.63638
$ws.Range("L83").Value = 54957.6
$ws.Range("M83").Value = 1561.36362
$ws.Range("N83").Value = -64941.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 24444
$ws.Range("I26").Value = 24444
$ws.Range("K26").Value = 24444
$ws.Range("M26").Value = -24114
$ws.Range("H112").Value = 20950
$ws.Range("J112").Value = 20950
$ws.Range("L112").Value = 20950
$ws.Range("N112").Value = -23904
$ws.Range("H122").Value = 3230.9
$ws.Range("I122").Value = 2954.9092
$ws.Range("J122").Value = 3568.2222
$ws.Range("K122").Value = 8864.7276
$ws.Range("L122").Value = 10704.6666
$ws.Range("M122").Value = -6414.7276
$ws.Range("N122").Value = -15604.6666
$ws.Range("H132").Value = 2352.5557
$ws.Range("I132").Value = 1821.8334
$ws.Range("K132").Value = 5465.5002
$ws.Range("M132").Value = -2935.5002

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H33").Value = 21663.334
$ws.Range("I33").Value = 5000
$ws.Range("J33").Value = 29995
$ws.Range("K33").Value = 5000
$ws.Range("L33").Value = 29995
$ws.Range("M33").Value = -4664
$ws.Range("N33").Value = -30667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3452397.5
$ws.Range("I31").Value = 5884767
$ws.Range("K31").Value = 5884767
$ws.Range("M31").Value = -5884472
$ws.Range("H32").Value = 3000
$ws.Range("I32").Value = 3000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -2684
$ws.Range("N32").ClearContents()
$ws.Range("H34").Value = 3452397.5
$ws.Range("I34").Value = 5884767
$ws.Range("K34").Value = 5884767
$ws.Range("M34").Value = -5884565
$ws.Range("H35").Value = 1214.2858
$ws.Range("I35").Value = 1214.2858
$ws.Range("K35").Value = 1214.2858
$ws.Range("M35").Value = -920.2858000000001
$ws.Range("H62").Value = 3809.1667
$ws.Range("I62").Value = 3252.5
$ws.Range("J62").Value = 4087.5
$ws.Range("K62").Value = 3252.5
$ws.Range("L62").Value = 4087.5
$ws.Range("M62").Value = -2628.5
$ws.Range("N62").Value = -5335.5
$ws.Range("H65").Value = 3809.1667
$ws.Range("I65").Value = 3252.5
$ws.Range("J65").Value = 4087.5
$ws.Range("K65").Value = 16262.5
$ws.Range("L65").Value = 20437.5
$ws.Range("M65").Value = -13142.5
$ws.Range("N65").Value = -26677.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 2984.25
$ws.Range("I94").Value = 1291.3334
$ws.Range("J94").Value = 4000
$ws.Range("K94").Value = 3874.0002
$ws.Range("L94").Value = 12000
$ws.Range("M94").Value = -3198.0002
$ws.Range("N94").Value = -13352
$ws.Range("H137").Value = 2489.1765
$ws.Range("I137").Value = 1922.2
$ws.Range("J137").Value = 2936.7896
$ws.Range("K137").Value = 5766.6
$ws.Range("L137").Value = 8810.3688
$ws.Range("M137").Value = -666.6000000000004
$ws.Range("N137").Value = -19010.3688

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 65881.31
$ws.Range("I102").Value = 2673.1428
$ws.Range("J102").Value = 115043.22
$ws.Range("K102").Value = 2673.1428
$ws.Range("L102").Value = 115043.22
$ws.Range("M102").Value = -1051.1428
$ws.Range("N102").Value = -118287.22
$ws.Range("H122").Value = 3944.257
$ws.Range("I122").Value = 3056.1875
$ws.Range("K122").Value = 9168.5625
$ws.Range("M122").Value = -6718.5625
$ws.Range("H139").Value = 33333.332
$ws.Range("J139").Value = 33333.332
$ws.Range("L139").Value = 33333.332
$ws.Range("N139").Value = -43613.332

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 125001720
$ws.Range("I22").Value = 200000350
$ws.Range("K22").Value = 200000350
$ws.Range("M22").Value = -200000055
$ws.Range("H27").Value = 125001720
$ws.Range("I27").Value = 200000350
$ws.Range("K27").Value = 200000350
$ws.Range("M27").Value = -200000243
$ws.Range("H32").Value = 29121.375
$ws.Range("I32").Value = 2585.2
$ws.Range("K32").Value = 2585.2
$ws.Range("M32").Value = -2268.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 772399.9
$ws.Range("I122").Value = 1669000.6
$ws.Range("J122").Value = 3884.8572
$ws.Range("K122").Value = 5007001.800000001
$ws.Range("L122").Value = 11654.5716
$ws.Range("M122").Value = -5004551.800000001
$ws.Range("N122").Value = -16554.5716
$ws.Range("H138").Value = 40000
$ws.Range("J138").Value = 40000
$ws.Range("L138").Value = 40000
$ws.Range("N138").Value = -50280

